{"js": "// Update the date line and the 25 division-problem answers in the table.\n// Each old value is unique in the document, so a simple search + replace\n// (matching the whole text of the run) is safe and keeps existing run\n// formatting (font, size, etc.) intact.\nconst replacements = [\n  [\"2024-11-05 Tuesday\", \"2024-11-06 Wednesday\"],\n  [\"414\u00f78=51, 6\", \"204\u00f72=102, 0\"],\n  [\"257\u00f76=42, 5\", \"485\u00f74=121, 1\"],\n  [\"434\u00f73=144, 2\", \"572\u00f78=71, 4\"],\n  [\"701\u00f76=116, 5\", \"494\u00f77=70, 4\"],\n  [\"628\u00f74=157, 0\", \"650\u00f72=325, 0\"],\n  [\"799\u00f75=159, 4\", \"820\u00f74=205, 0\"],\n  [\"537\u00f73=179, 0\", \"800\u00f79=88, 8\"],\n  [\"105\u00f79=11, 6\", \"441\u00f77=63, 0\"],\n  [\"476\u00f74=119, 0\", \"308\u00f74=77, 0\"],\n  [\"879\u00f76=146, 3\", \"854\u00f72=427, 0\"],\n  [\"960\u00f72=480, 0\", \"221\u00f79=24, 5\"],\n  [\"502\u00f77=71, 5\", \"407\u00f77=58, 1\"],\n  [\"811\u00f79=90, 1\", \"886\u00f79=98, 4\"],\n  [\"770\u00f78=96, 2\", \"641\u00f74=160, 1\"],\n  [\"733\u00f73=244, 1\", \"984\u00f75=196, 4\"],\n  [\"835\u00f76=139, 1\", \"360\u00f72=180, 0\"],\n  [\"346\u00f79=38, 4\", \"230\u00f75=46, 0\"],\n  [\"911\u00f77=130, 1\", \"650\u00f77=92, 6\"],\n  [\"726\u00f75=145, 1\", \"900\u00f73=300, 0\"],\n  [\"238\u00f74=59, 2\", \"921\u00f74=230, 1\"],\n  [\"365\u00f72=182, 1\", \"276\u00f73=92, 0\"],\n  [\"600\u00f78=75, 0\", \"818\u00f78=102, 2\"],\n  [\"647\u00f78=80, 7\", \"942\u00f77=134, 4\"],\n  [\"789\u00f79=87, 6\", \"742\u00f78=92, 6\"],\n  [\"137\u00f78=17, 1\", \"913\u00f79=101, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 division-problem answers in the table.\n# Each old value is unique in the document, so Find/Replace on the whole\n# text keeps the existing run formatting (font, size, etc.) intact.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-11-05 Tuesday\", \"2024-11-06 Wednesday\"),\n    @(\"414\u00f78=51, 6\", \"204\u00f72=102, 0\"),\n    @(\"257\u00f76=42, 5\", \"485\u00f74=121, 1\"),\n    @(\"434\u00f73=144, 2\", \"572\u00f78=71, 4\"),\n    @(\"701\u00f76=116, 5\", \"494\u00f77=70, 4\"),\n    @(\"628\u00f74=157, 0\", \"650\u00f72=325, 0\"),\n    @(\"799\u00f75=159, 4\", \"820\u00f74=205, 0\"),\n    @(\"537\u00f73=179, 0\", \"800\u00f79=88, 8\"),\n    @(\"105\u00f79=11, 6\", \"441\u00f77=63, 0\"),\n    @(\"476\u00f74=119, 0\", \"308\u00f74=77, 0\"),\n    @(\"879\u00f76=146, 3\", \"854\u00f72=427, 0\"),\n    @(\"960\u00f72=480, 0\", \"221\u00f79=24, 5\"),\n    @(\"502\u00f77=71, 5\", \"407\u00f77=58, 1\"),\n    @(\"811\u00f79=90, 1\", \"886\u00f79=98, 4\"),\n    @(\"770\u00f78=96, 2\", \"641\u00f74=160, 1\"),\n    @(\"733\u00f73=244, 1\", \"984\u00f75=196, 4\"),\n    @(\"835\u00f76=139, 1\", \"360\u00f72=180, 0\"),\n    @(\"346\u00f79=38, 4\", \"230\u00f75=46, 0\"),\n    @(\"911\u00f77=130, 1\", \"650\u00f77=92, 6\"),\n    @(\"726\u00f75=145, 1\", \"900\u00f73=300, 0\"),\n    @(\"238\u00f74=59, 2\", \"921\u00f74=230, 1\"),\n    @(\"365\u00f72=182, 1\", \"276\u00f73=92, 0\"),\n    @(\"600\u00f78=75, 0\", \"818\u00f78=102, 2\"),\n    @(\"647\u00f78=80, 7\", \"942\u00f77=134, 4\"),\n    @(\"789\u00f79=87, 6\", \"742\u00f78=92, 6\"),\n    @(\"137\u00f78=17, 1\", \"913\u00f79=101, 4\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
